# New weekly price record for "Poroto verde" (Terminal Hortofrutícola Agro
# Chillán) is inserted at row 11, pushing all existing historical rows
# (11-78) down by one (to 12-79).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11 (shifts rows 11..78 down to 12..79,
# inheriting the formatting of the row above - same as Excel's UI "Insert").
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A11").Value = 7
$ws.Range("B11").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C11").Value = "Ñuble"
$ws.Range("D11").Value = 44635
$ws.Range("E11").Value = 16
$ws.Range("F11").Value = 100112031
$ws.Range("G11").Value = "Poroto verde"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 80
$ws.Range("K11").Value = 29000
$ws.Range("L11").Value = 30000
$ws.Range("M11").Value = 29500
$ws.Range("N11").Value = "$/saco 25 kilos"
$ws.Range("O11").Value = "Región del Maule"
$ws.Range("P11").Value = 1180
$ws.Range("Q11").Value = 25
$ws.Range("R11").Value = "Hortaliza"
